$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.313179
$ws.Range("H2").Value = 0.939537
$ws.Range("I2").Value = 0.02707464596575709
$ws.Range("J2").Value = 0.0270746459657571
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3659943333333333
$ws.Range("N2").Value = 1.097983
$ws.Range("O2").Value = 0.006726051721149161
$ws.Range("P2").Value = 0.006726051721149162
$ws.Range("Q2").Value = 0.114621739319
$ws.Range("R2").Value = 1.031595653871
$ws.Range("S2").Value = 0.0001821054690974847
$ws.Range("T2").Value = 0.0001821054690974847

# Row 3
$ws.Range("G3").Value = 0.313179
$ws.Range("H3").Value = 0.939537
$ws.Range("I3").Value = 0.02707464596575709
$ws.Range("J3").Value = 0.0270746459657571
$ws.Range("O3").Value = 0.001678071748088335
$ws.Range("P3").Value = 0.001678071748088335
$ws.Range("Q3").Value = 0.028596792062
$ws.Range("R3").Value = 0.257371128558
$ws.Range("S3").Value = 0.0000454331984846308
$ws.Range("T3").Value = 0.00004543319848463081

# Row 4
$ws.Range("G4").Value = 0.313179
$ws.Range("H4").Value = 0.939537
$ws.Range("I4").Value = 0.02707464596575709
$ws.Range("J4").Value = 0.0270746459657571
$ws.Range("M4").Value = 53.897087
$ws.Range("N4").Value = 161.691261
$ws.Range("O4").Value = 0.9904923704135933
$ws.Range("P4").Value = 0.9904923704135934
$ws.Range("Q4").Value = 16.879435809573
$ws.Range("R4").Value = 151.914922286157
$ws.Range("S4").Value = 0.02681723026073157
$ws.Range("T4").Value = 0.02681723026073158

# Row 5
$ws.Range("G5").Value = 0.313179
$ws.Range("H5").Value = 0.939537
$ws.Range("I5").Value = 0.02707464596575709
$ws.Range("J5").Value = 0.0270746459657571
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.06004666666666667
$ws.Range("N5").Value = 0.18014
$ws.Range("O5").Value = 0.001103506117169219
$ws.Range("P5").Value = 0.001103506117169219
$ws.Range("Q5").Value = 0.01880535502
$ws.Range("R5").Value = 0.16924819518
$ws.Range("S5").Value = 0.00002987703744340386
$ws.Range("T5").Value = 0.00002987703744340387

# Row 6
$ws.Range("I6").Value = 0.6982806158817221
$ws.Range("J6").Value = 0.6982806158817222
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.3659943333333333
$ws.Range("N6").Value = 1.097983
$ws.Range("O6").Value = 0.006726051721149161
$ws.Range("P6").Value = 0.006726051721149162
$ws.Range("Q6").Value = 2.956202597305778
$ws.Range("R6").Value = 26.605823375752
$ws.Range("S6").Value = 0.004696671538296353
$ws.Range("T6").Value = 0.004696671538296355

# Row 7
$ws.Range("I7").Value = 0.6982806158817221
$ws.Range("J7").Value = 0.6982806158817222
$ws.Range("O7").Value = 0.001678071748088335
$ws.Range("P7").Value = 0.001678071748088335
$ws.Range("S7").Value = 0.001171764973748841
$ws.Range("T7").Value = 0.001171764973748841

# Row 8
$ws.Range("I8").Value = 0.6982806158817221
$ws.Range("J8").Value = 0.6982806158817222
$ws.Range("M8").Value = 53.897087
$ws.Range("N8").Value = 161.691261
$ws.Range("O8").Value = 0.9904923704135933
$ws.Range("P8").Value = 0.9904923704135934
$ws.Range("Q8").Value = 435.3365450374427
$ws.Range("R8").Value = 3918.028905336984
$ws.Range("S8").Value = 0.6916416224385508
$ws.Range("T8").Value = 0.6916416224385509

# Row 9
$ws.Range("I9").Value = 0.6982806158817221
$ws.Range("J9").Value = 0.6982806158817222
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.06004666666666667
$ws.Range("N9").Value = 0.18014
$ws.Range("O9").Value = 0.001103506117169219
$ws.Range("P9").Value = 0.001103506117169219
$ws.Range("Q9").Value = 0.485007815128889
$ws.Range("R9").Value = 4.365070336160001
$ws.Range("S9").Value = 0.0007705569311261697
$ws.Range("T9").Value = 0.0007705569311261699

# Row 10
$ws.Range("G10").Value = 2.897745666666667
$ws.Range("H10").Value = 8.693237
$ws.Range("I10").Value = 0.2505130868410934
$ws.Range("J10").Value = 0.2505130868410934
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.3659943333333333
$ws.Range("N10").Value = 1.097983
$ws.Range("O10").Value = 0.006726051721149161
$ws.Range("P10").Value = 0.006726051721149162
$ws.Range("Q10").Value = 1.060558493441222
$ws.Range("R10").Value = 9.545026440970998
$ws.Range("S10").Value = 0.001684963978917925
$ws.Range("T10").Value = 0.001684963978917925

# Row 11
$ws.Range("G11").Value = 2.897745666666667
$ws.Range("H11").Value = 8.693237
$ws.Range("I11").Value = 0.2505130868410934
$ws.Range("J11").Value = 0.2505130868410934
$ws.Range("O11").Value = 0.001678071748088335
$ws.Range("P11").Value = 0.001678071748088335
$ws.Range("Q11").Value = 0.2645970204842222
$ws.Range("R11").Value = 2.381373184358
$ws.Range("S11").Value = 0.0004203789335544385
$ws.Range("T11").Value = 0.0004203789335544385

# Row 12
$ws.Range("G12").Value = 2.897745666666667
$ws.Range("H12").Value = 8.693237
$ws.Range("I12").Value = 0.2505130868410934
$ws.Range("J12").Value = 0.2505130868410934
$ws.Range("M12").Value = 53.897087
$ws.Range("N12").Value = 161.691261
$ws.Range("O12").Value = 0.9904923704135933
$ws.Range("P12").Value = 0.9904923704135934
$ws.Range("Q12").Value = 156.1800503002063
$ws.Range("R12").Value = 1405.620452701857
$ws.Range("S12").Value = 0.2481313012048609
$ws.Range("T12").Value = 0.2481313012048609

# Row 13
$ws.Range("G13").Value = 2.897745666666667
$ws.Range("H13").Value = 8.693237
$ws.Range("I13").Value = 0.2505130868410934
$ws.Range("J13").Value = 0.2505130868410934
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.06004666666666667
$ws.Range("N13").Value = 0.18014
$ws.Range("O13").Value = 0.001103506117169219
$ws.Range("P13").Value = 0.001103506117169219
$ws.Range("Q13").Value = 0.1739999681311111
$ws.Range("R13").Value = 1.56599971318
$ws.Range("S13").Value = 0.0002764427237600902
$ws.Range("T13").Value = 0.0002764427237600902

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.2791366666666666
$ws.Range("H14").Value = 0.83741
$ws.Range("I14").Value = 0.02413165131142748
$ws.Range("J14").Value = 0.02413165131142749
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.3659943333333333
$ws.Range("N14").Value = 1.097983
$ws.Range("O14").Value = 0.006726051721149161
$ws.Range("P14").Value = 0.006726051721149162
$ws.Range("Q14").Value = 0.1021624382255555
$ws.Range("R14").Value = 0.9194619440299999
$ws.Range("S14").Value = 0.0001623107348373982
$ws.Range("T14").Value = 0.0001623107348373983

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.2791366666666666
$ws.Range("H15").Value = 0.83741
$ws.Range("I15").Value = 0.02413165131142748
$ws.Range("J15").Value = 0.02413165131142749
$ws.Range("O15").Value = 0.001678071748088335
$ws.Range("P15").Value = 0.001678071748088335
$ws.Range("Q15").Value = 0.02548834121555556
$ws.Range("R15").Value = 0.22939507094
$ws.Range("S15").Value = 0.00004049464230042529
$ws.Range("T15").Value = 0.00004049464230042531

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.2791366666666666
$ws.Range("H16").Value = 0.83741
$ws.Range("I16").Value = 0.02413165131142748
$ws.Range("J16").Value = 0.02413165131142749
$ws.Range("M16").Value = 53.897087
$ws.Range("N16").Value = 161.691261
$ws.Range("O16").Value = 0.9904923704135933
$ws.Range("P16").Value = 0.9904923704135934
$ws.Range("Q16").Value = 15.04465320822333
$ws.Range("R16").Value = 135.40187887401
$ws.Range("S16").Value = 0.02390221650945011
$ws.Range("T16").Value = 0.02390221650945011

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.2791366666666666
$ws.Range("H17").Value = 0.83741
$ws.Range("I17").Value = 0.02413165131142748
$ws.Range("J17").Value = 0.02413165131142749
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.06004666666666667
$ws.Range("N17").Value = 0.18014
$ws.Range("O17").Value = 0.001103506117169219
$ws.Range("P17").Value = 0.001103506117169219
$ws.Range("Q17").Value = 0.01676122637777778
$ws.Range("R17").Value = 0.1508510374
$ws.Range("S17").Value = 0.00002662942483955484
$ws.Range("T17").Value = 0.00002662942483955484
